$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Junio de 2020 a las 21:54"

# Row 4
$ws.Cells.Item(4, 2).Value = 2175058
$ws.Cells.Item(4, 3).Value = 12830
$ws.Cells.Item(4, 4).Value = 875189
$ws.Cells.Item(4, 5).Value = 1181736
$ws.Cells.Item(4, 7).Value = 275
$ws.Cells.Item(4, 8).Value = 118133

# Row 13
$ws.Cells.Item(13, 2).Value = 187967
$ws.Cells.Item(13, 3).Value = 296
$ws.Cells.Item(13, 5).Value = 6490
$ws.Cells.Item(13, 7).Value = 7
$ws.Cells.Item(13, 8).Value = 8877

# Row 16
$ws.Cells.Item(16, 2).Value = 157372
$ws.Cells.Item(16, 3).Value = 152
$ws.Cells.Item(16, 4).Value = 73044
$ws.Cells.Item(16, 5).Value = 54892
$ws.Cells.Item(16, 7).Value = 29
$ws.Cells.Item(16, 8).Value = 29436

# Row 20
$ws.Cells.Item(20, 2).Value = 99073
$ws.Cells.Item(20, 3).Value = 286
$ws.Cells.Item(20, 4).Value = 60999
$ws.Cells.Item(20, 5).Value = 29900

# Row 31
$ws.Cells.Item(31, 2).Value = 46289
$ws.Cells.Item(31, 3).Value = 1691
$ws.Cells.Item(31, 4).Value = 12329
$ws.Cells.Item(31, 5).Value = 32288
$ws.Cells.Item(31, 7).Value = 97
$ws.Cells.Item(31, 8).Value = 1672

# Row 75
$ws.Cells.Item(75, 1).Value = "Costa de Marfil"
$ws.Cells.Item(75, 2).Value = 5439
$ws.Cells.Item(75, 3).Value = 355
$ws.Cells.Item(75, 4).Value = 2590
$ws.Cells.Item(75, 5).Value = 2803
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 46

# Row 76
$ws.Cells.Item(76, 1).Value = "Uzbekistan"
$ws.Cells.Item(76, 2).Value = 5263
$ws.Cells.Item(76, 3).Value = 183
$ws.Cells.Item(76, 4).Value = 4019
$ws.Cells.Item(76, 5).Value = 1225
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 19

# Row 77
$ws.Cells.Item(77, 1).Value = "Senegal"
$ws.Cells.Item(77, 2).Value = 5173
$ws.Cells.Item(77, 3).Value = 83
$ws.Cells.Item(77, 4).Value = 3424
$ws.Cells.Item(77, 5).Value = 1685
$ws.Cells.Item(77, 7).Value = 4
$ws.Cells.Item(77, 8).Value = 64

# Row 78
$ws.Cells.Item(78, 1).Value = "Tayikistan"
$ws.Cells.Item(78, 2).Value = 5097
$ws.Cells.Item(78, 3).Value = 62
$ws.Cells.Item(78, 4).Value = 3503
$ws.Cells.Item(78, 5).Value = 1544
$ws.Cells.Item(78, 8).Value = 50

# Row 108
$ws.Cells.Item(108, 2).Value = 1744
$ws.Cells.Item(108, 3).Value = 29
$ws.Cells.Item(108, 4).Value = 771
$ws.Cells.Item(108, 5).Value = 961

# Row 130
$ws.Cells.Item(130, 1).Value = "Congo"
$ws.Cells.Item(130, 2).Value = 883
$ws.Cells.Item(130, 3).Value = 155
$ws.Cells.Item(130, 4).Value = 391
$ws.Cells.Item(130, 5).Value = 465
$ws.Cells.Item(130, 7).Value = 3
$ws.Cells.Item(130, 8).Value = 27

# Row 131
$ws.Cells.Item(131, 1).Value = "Georgia"
$ws.Cells.Item(131, 2).Value = 879
$ws.Cells.Item(131, 3).Value = 15
$ws.Cells.Item(131, 4).Value = 704
$ws.Cells.Item(131, 5).Value = 161
$ws.Cells.Item(131, 8).Value = 14

# Row 132
$ws.Cells.Item(132, 1).Value = "Principado de Andorra"
$ws.Cells.Item(132, 2).Value = 853
$ws.Cells.Item(132, 4).Value = 789
$ws.Cells.Item(132, 5).Value = 13
$ws.Cells.Item(132, 8).Value = 51

# Row 133
$ws.Cells.Item(133, 1).Value = "Republica del Chad"
$ws.Cells.Item(133, 2).Value = 850
$ws.Cells.Item(133, 4).Value = 720
$ws.Cells.Item(133, 5).Value = 57
$ws.Cells.Item(133, 8).Value = 73

# Row 134
$ws.Cells.Item(134, 1).Value = "Uruguay"
$ws.Cells.Item(134, 2).Value = 848
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = 788
$ws.Cells.Item(134, 5).Value = 37
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 23

# Row 135
$ws.Cells.Item(135, 1).Value = "Yemen"
$ws.Cells.Item(135, 2).Value = 844
$ws.Cells.Item(135, 3).Value = 116
$ws.Cells.Item(135, 4).Value = 79
$ws.Cells.Item(135, 5).Value = 557
$ws.Cells.Item(135, 7).Value = 44
$ws.Cells.Item(135, 8).Value = 208

# Row 136
$ws.Cells.Item(136, 1).Value = "Cabo Verde"
$ws.Cells.Item(136, 2).Value = 759
$ws.Cells.Item(136, 3).Value = 9
$ws.Cells.Item(136, 4).Value = 301
$ws.Cells.Item(136, 5).Value = 452
$ws.Cells.Item(136, 8).Value = 6

# Row 143
$ws.Cells.Item(143, 1).Value = "Ruanda"
$ws.Cells.Item(143, 2).Value = 612
$ws.Cells.Item(143, 3).Value = 30
$ws.Cells.Item(143, 4).Value = 338
$ws.Cells.Item(143, 5).Value = 272
$ws.Cells.Item(143, 8).Value = 2

# Row 144
$ws.Cells.Item(144, 1).Value = "Mozambique"
$ws.Cells.Item(144, 2).Value = 609
$ws.Cells.Item(144, 3).Value = 26
$ws.Cells.Item(144, 4).Value = 157
$ws.Cells.Item(144, 5).Value = 449
$ws.Cells.Item(144, 8).Value = 3

# Row 148
$ws.Cells.Item(148, 1).Value = "Suazilandia"
$ws.Cells.Item(148, 2).Value = 506
$ws.Cells.Item(148, 3).Value = 16
$ws.Cells.Item(148, 4).Value = 249
$ws.Cells.Item(148, 5).Value = 253
$ws.Cells.Item(148, 8).Value = 4

# Row 149
$ws.Cells.Item(149, 1).Value = "Estado de Palestina"
$ws.Cells.Item(149, 2).Value = 505
$ws.Cells.Item(149, 3).Value = 13
$ws.Cells.Item(149, 4).Value = 415
$ws.Cells.Item(149, 5).Value = 87
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 3

# Row 150
$ws.Cells.Item(150, 1).Value = "Liberia"
$ws.Cells.Item(150, 2).Value = 498
$ws.Cells.Item(150, 3).Value = 40
$ws.Cells.Item(150, 4).Value = 221
$ws.Cells.Item(150, 5).Value = 244
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = 33

# Row 151
$ws.Cells.Item(151, 1).Value = "Reunion"
$ws.Cells.Item(151, 2).Value = 496
$ws.Cells.Item(151, 3).Value = 1
$ws.Cells.Item(151, 4).Value = 460
$ws.Cells.Item(151, 5).Value = 35
$ws.Cells.Item(151, 8).Value = 1

# Row 159
$ws.Cells.Item(159, 2).Value = 326
$ws.Cells.Item(159, 3).Value = 1
$ws.Cells.Item(159, 5).Value = 2

# Row 166
$ws.Cells.Item(166, 4).Value = 78
$ws.Cells.Item(166, 5).Value = 93

# Row 171
$ws.Cells.Item(171, 2).Value = 144
$ws.Cells.Item(171, 3).Value = 2
$ws.Cells.Item(171, 5).Value = 8

# Row 172
$ws.Cells.Item(172, 1).Value = "Angola"
$ws.Cells.Item(172, 2).Value = 142
$ws.Cells.Item(172, 3).Value = 2
$ws.Cells.Item(172, 4).Value = 64
$ws.Cells.Item(172, 5).Value = 72
$ws.Cells.Item(172, 8).Value = 6

# Row 173
$ws.Cells.Item(173, 1).Value = "Brunei"
$ws.Cells.Item(173, 2).Value = 141
$ws.Cells.Item(173, 4).Value = 138
$ws.Cells.Item(173, 5).Value = 1
$ws.Cells.Item(173, 8).Value = 2

# Row 192
$ws.Cells.Item(192, 2).Value = 30
$ws.Cells.Item(192, 3).Value = 2
$ws.Cells.Item(192, 5).Value = 5

# Row 208
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 8).Value = 0

# Row 209
$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 1

# Row 210
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

# Row 211
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1

# Row 213
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

# Row 214
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
